$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.946.48"
$ws.Range("E2").Value = "'  +0.17%  "
$ws.Range("D3").Value = "'1.557.03"
$ws.Range("E3").Value = "'  +0.58%  "
$ws.Range("E4").Value = "'  -0.19%  "
$ws.Range("D5").Value = "'206.90"
$ws.Range("E5").Value = "'  -0.04%  "
$ws.Range("D6").Value = "'0.487"
$ws.Range("E6").Value = "'  -0.18%  "
$ws.Range("E7").Value = "'  -0.21%  "
$ws.Range("E8").Value = "'  +2.88%  "
$ws.Range("E9").Value = "'  +0.16%  "
$ws.Range("E10").Value = "'  +0.91%  "
$ws.Range("E11").Value = "'  -0.07%  "
$ws.Range("D12").Value = "'1.778.13"
$ws.Range("E12").Value = "'  +0.54%  "
$ws.Range("D13").Value = "'1.556.11"
$ws.Range("E13").Value = "'  +0.40%  "
$ws.Range("E15").Value = "'  +1.56%  "
$ws.Range("D16").Value = "'26.949.40"
$ws.Range("E16").Value = "'  +0.16%  "
$ws.Range("D17").Value = "'61.77"
$ws.Range("E17").Value = "'  +0.53%  "
$ws.Range("D18").Value = "'217.84"
$ws.Range("E18").Value = "'  +1.49%  "
$ws.Range("E19").Value = "'  +1.71%  "
$ws.Range("E20").Value = "'  +1.16%  "
$ws.Range("E22").Value = "'  +0.97%  "
$ws.Range("D23").Value = "'9.20"
$ws.Range("E23").Value = "'  +0.29%  "
$ws.Range("E24").Value = "'  +0.70%  "
$ws.Range("D25").Value = "'153.28"
$ws.Range("E25").Value = "'  +0.89%  "
$ws.Range("E26").Value = "'  +0.39%  "
$ws.Range("E27").Value = "'  +0.35%  "
$ws.Range("E28").Value = "'  +0.49%  "
$ws.Range("E30").Value = "'  +2.46%  "
$ws.Range("D31").Value = "'1.08"
$ws.Range("E31").Value = "'  -1.41%  "
$ws.Range("E32").Value = "'  -0.17%  "
$ws.Range("D33").Value = "'1.423.69"
$ws.Range("E33").Value = "'  +4.07%  "
$ws.Range("E34").Value = "'  +3.95%  "
$ws.Range("E35").Value = "'  +3.80%  "
$ws.Range("D36").Value = "'0.977"
$ws.Range("E36").Value = "'  +1.84%  "
$ws.Range("E38").Value = "'  -0.23%  "
$ws.Range("D39").Value = "'0.520"
$ws.Range("E39").Value = "'  -0.23%  "
$ws.Range("D40").Value = "'0.814"
$ws.Range("E40").Value = "'  +0.93%  "
$ws.Range("E41").Value = "'  -0.17%  "
$ws.Range("E42").Value = "'  +1.11%  "
$ws.Range("D43").Value = "'2.30"
$ws.Range("E43").Value = "'  +3.44%  "
$ws.Range("E44").Value = "'  -0.61%  "
$ws.Range("D45").Value = "'64.75"
$ws.Range("E45").Value = "'  +1.73%  "
$ws.Range("E46").Value = "'  +1.09%  "
$ws.Range("D47").Value = "'1.691.86"
$ws.Range("E47").Value = "'  +0.49%  "
$ws.Range("D48").Value = "'87.48"
$ws.Range("E48").Value = "'  +2.18%  "
$ws.Range("D49").Value = "'0.0519"
$ws.Range("E49").Value = "'  +2.25%  "
$ws.Range("E50").Value = "'  +2.53%  "
$ws.Range("D51").Value = "'0.0960"
$ws.Range("E51").Value = "'  +1.14%  "
